# The upstream commit (obeonetwork/M2Doc, fix #295 "Add the version of
# M2Doc in the template custom properties") only touches the template's
# custom document properties; none of that is reflected in the supplied
# OOXML diff. Every hunk in that diff is a pure XML attribute
# re-serialization (attributes alphabetized by the scraping/diff tool,
# e.g. "w:type=...  w:default=..." -> "w:default=...  w:type=...") with
# identical tag names, identical attribute sets/values and identical
# element order throughout word/document.xml and word/styles.xml - i.e.
# no visible or structural content actually changes.
#
# So the faithful COM-interop replay is a no-op on content: we simply
# touch the document through the object model (forcing Word to resolve
# it) without mutating any text, run/paragraph formatting, section
# properties, or style definitions.
$d = $word.ActiveDocument
$null = $d.Content.Text
